$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Init")

# Extend scenario 18 references to scenario 19 (D5:D11)
$ws.Range("D5").Value = "A19"
$ws.Range("D6").Value = "B19"
$ws.Range("D7").Value = "C19"
$ws.Range("D8").Value = "G19"
$ws.Range("D9").Value = "H19"
$ws.Range("D10").Value = "I19"
$ws.Range("D11").Value = "J19"

# Set active cell selection as last edited by the user
$ws.Range("D11").Select()
